# "Slight Changes in Work Schedule"
#
# - Update three "Completed on" dates from 2018-06-05 (43256) to
#   2018-09-05 (43348): D3, C4, C5 on Sheet1.
# - Bold the header row (A1:D1).
# - Set the sheet's page orientation to portrait.
# - Move the active selection from A6 to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Completed on" dates ---
$ws.Range("D3").Value = 43348
$ws.Range("C4").Value = 43348
$ws.Range("C5").Value = 43348

# --- Bold the header row ---
$ws.Range("A1:D1").Font.Bold = $true

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Move selection to D6 (matches the saved cursor position) ---
$ws.Range("D6").Select()
